$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.005" are not
# reinterpreted/rounded as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.924.96"
$ws.Range("E2").Value = "  -0.60%  "

$ws.Range("D3").Value = "1.643.63"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.64%  "

$ws.Range("D5").Value = "215.69"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").Value = "0.5058"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").Value = "1.006"
$ws.Range("E7").Value = "  -0.61%  "

$ws.Range("D8").Value = "0.2574"
$ws.Range("E8").Value = "  -0.42%  "

$ws.Range("D9").Value = "0.06405"
$ws.Range("E9").Value = "  -0.60%  "

$ws.Range("D10").Value = "19.64"
$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("D11").Value = "0.07812"
$ws.Range("E11").Value = "  +0.81%  "

$ws.Range("D12").Value = "1.658.78"
$ws.Range("E12").Value = "  +0.75%  "

$ws.Range("D13").Value = "4.284"
$ws.Range("E13").Value = "  +0.50%  "

$ws.Range("D14").Value = "1.867.67"
$ws.Range("E14").Value = "  -0.28%  "

$ws.Range("D15").Value = "0.5440"
$ws.Range("E15").Value = "  -0.38%  "

$ws.Range("D16").Value = "0.0₅7881"
$ws.Range("E16").Value = "  -0.66%  "

$ws.Range("D17").Value = "64.89"
$ws.Range("E17").Value = "  +1.77%  "

$ws.Range("D18").Value = "25.977.46"
$ws.Range("E18").Value = "  -0.41%  "

$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  -0.56%  "

$ws.Range("D20").Value = "198.55"
$ws.Range("E20").Value = "  -2.70%  "

$ws.Range("D21").Value = "4.398"
$ws.Range("E21").Value = "  +1.95%  "

$ws.Range("D22").Value = "9.988"
$ws.Range("E22").Value = "  -0.28%  "

$ws.Range("D23").Value = "5.985"
$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("D24").Value = "1.008"
$ws.Range("E24").Value = "  -0.52%  "

$ws.Range("D25").Value = "1.876"
$ws.Range("E25").Value = "  -3.26%  "

$ws.Range("D26").Value = "140.34"
$ws.Range("E26").Value = "  -1.21%  "

$ws.Range("D27").Value = "0.1146"
$ws.Range("E27").Value = "  -1.08%  "

$ws.Range("D28").Value = "6.866"
$ws.Range("E28").Value = "  +1.71%  "

$ws.Range("D29").Value = "15.74"
$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("D30").Value = "1.246"
$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").Value = "0.04921"
$ws.Range("E31").Value = "  -2.86%  "

$ws.Range("D32").Value = "3.270"
$ws.Range("E32").Value = "  +0.18%  "

$ws.Range("D33").Value = "3.202"
$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("D34").Value = "1.538"
$ws.Range("E34").Value = "  -0.55%  "

$ws.Range("D35").Value = "2.379"
$ws.Range("E35").Value = "  +1.29%  "

$ws.Range("D36").Value = "0.8956"
$ws.Range("E36").Value = "  -0.26%  "

$ws.Range("D37").Value = "2.606"
$ws.Range("E37").Value = "  -0.58%  "

$ws.Range("D38").Value = "1.141.25"
$ws.Range("E38").Value = "  -0.94%  "

$ws.Range("D39").Value = "0.5560"
$ws.Range("E39").Value = "  -1.49%  "

$ws.Range("D40").Value = "0.01564"
$ws.Range("E40").Value = "  -0.68%  "

$ws.Range("D41").Value = "1.009"
$ws.Range("E41").Value = "  -0.32%  "

$ws.Range("D42").Value = "5.703"
$ws.Range("E42").Value = "  +0.57%  "

$ws.Range("D43").Value = "0.8201"
$ws.Range("E43").Value = "  +0.48%  "

$ws.Range("D44").Value = "99.58"
$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("D45").Value = "0.0₈120"
$ws.Range("E45").Value = "  +6.60%  "

$ws.Range("D46").Value = "1.778.21"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("D47").Value = "0.4525"
$ws.Range("E47").Value = "  -0.36%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "55.47"
$ws.Range("E48").Value = "  +0.81%  "

$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").Value = "1.005"
$ws.Range("E49").Value = "  -0.61%  "

$ws.Range("D50").Value = "0.05088"
$ws.Range("E50").Value = "  +0.81%  "

$ws.Range("E51").Value = "  -0.49%  "
